$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Rewrite the "Além disso, ..." sentence: reorder/reword the clause
#    and recolor part of it red. This also removes the old proofErr
#    markers and the old _GoBack bookmark that used to sit inside this
#    run of text (the replace below swallows that span).
# ---------------------------------------------------------------------
$old = "rotinas específica em linguagem C foram desenvolvidas  Este código"
$new = "foram desenvolvidas rotinas específicas em linguagem C. Este código"
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# Color the five red phrases individually so each keeps its own run
# (matches the five separate <w:r> elements in the target markup).
$seg = $d.Content
$seg.Find.Execute("Além disso, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg.Font.Color = 255

$seg = $d.Content
$seg.Find.Execute("foram desenvolvidas ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg.Font.Color = 255

$seg = $d.Content
$seg.Find.Execute("rotinas específicas", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg.Font.Color = 255

$seg = $d.Content
$seg.Find.Execute(" em linguagem ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg.Font.Color = 255

$seg = $d.Content
$seg.Find.Execute("C. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seg.Font.Color = 255

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark so it now sits right after the
#    "Código numérico para solução das análises" paragraph's text.
#    (Adding a bookmark named "_GoBack" elsewhere automatically moves
#    it, since bookmark names are unique within the document, and the
#    old one was already consumed by the replace above.)
#
#    A collapsed range sitting exactly on that paragraph's end-of-text
#    boundary is not accepted directly, so nudge it into a reachable
#    spot by temporarily inserting a marker, bookmarking before it,
#    then removing the marker again.
# ---------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("Código numérico para solução das análises", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$afterPoint = $d.Range($target.End, $target.End)
$afterPoint.InsertAfter("ZZBOOKMARKMARKERZZ")

$marker = $d.Content
$marker.Find.Execute("ZZBOOKMARKMARKERZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$bmSpot = $d.Range($marker.Start, $marker.Start)
$bmSpot.Bookmarks.Add("_GoBack")

$marker = $d.Content
$marker.Find.Execute("ZZBOOKMARKMARKERZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$marker.Text = ""
